# RDCC-5182 fix for version compare in SRD file
# Adds a new "VERSION" worksheet as the first sheet in the workbook, used
# to hold the file/template version so it can be compared during upload.

$wb = $excel.ActiveWorkbook

# Worksheets.Add() inserts the new sheet before the active sheet, which in
# this workbook puts it first (matches the target: VERSION becomes sheet
# index 1, ahead of "Staff Data" and "Sheet2").
$versionSheet = $wb.Worksheets.Add()
$versionSheet.Name = "VERSION"

# Populate the version marker cells.
$versionSheet.Range("A6").Value = "File version"
$versionSheet.Range("B6").Value = "vx.xx"

# Match the author's saved selection/active cell on the new sheet.
$versionSheet.Range("B6").Select() | Out-Null
